$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are numeric-looking text (e.g. thousand-dot-separated
# prices, or decimals with significant trailing zeros) that must stay TEXT, so
# force the text number format before assigning to avoid Excel auto-converting
# them to actual numbers (which would also mangle values like 61.764.57).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.764.57'
$ws.Range('E2').Value = '  -8.11%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.907.44'
$ws.Range('E3').Value = '  -10.32%  '

$ws.Range('E4').Value = '  +0.36%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '518.26'
$ws.Range('E5').Value = '  -13.01%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '123.05'
$ws.Range('E6').Value = '  -19.35%  '

$ws.Range('E7').Value = '  +0.01%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '2.898.11'
$ws.Range('E8').Value = '  -10.31%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.442'
$ws.Range('E9').Value = '  -18.67%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.138'
$ws.Range('E10').Value = '  -19.89%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.52'
$ws.Range('E11').Value = '  -14.51%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.409'
$ws.Range('E12').Value = '  -16.92%  '

$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '29.78'
$ws.Range('E13').Value = '  -23.66%  '

$ws.Range('B14').Value = 'ShibaInu'
$ws.Range('C14').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000191'
$ws.Range('E14').Value = '  -21.02%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.388.52'
$ws.Range('E15').Value = '  -10.08%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '61.775.95'
$ws.Range('E16').Value = '  -8.14%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.109'
$ws.Range('E17').Value = '  -4.23%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.926.89'
$ws.Range('E18').Value = '  -9.72%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '456.85'
$ws.Range('E19').Value = '  -14.08%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.82'
$ws.Range('E20').Value = '  -17.57%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.05'
$ws.Range('E21').Value = '  -18.89%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.608'
$ws.Range('E22').Value = '  -19.78%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.15'
$ws.Range('E23').Value = '  -22.20%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '72.86'
$ws.Range('E24').Value = '  -14.66%  '

$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.994'
$ws.Range('E25').Value = '  -0.45%  '

$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.20'
$ws.Range('E26').Value = '  -18.31%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.56'
$ws.Range('E27').Value = '  -19.85%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.53'
$ws.Range('E28').Value = '  -18.35%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '23.67'
$ws.Range('E29').Value = '  -18.62%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.72'
$ws.Range('E30').Value = '  -21.19%  '

$ws.Range('B31').Value = 'FirstDigitalUSD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.01'
$ws.Range('E31').Value = '  +0.57%  '

$ws.Range('B32').Value = 'Mantle'
$ws.Range('C32').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.04'
$ws.Range('E32').Value = '  -9.70%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.16'
$ws.Range('E33').Value = '  -19.19%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '50.74'
$ws.Range('E34').Value = '  -5.33%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '440.79'
$ws.Range('E35').Value = '  -17.09%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.22'
$ws.Range('E36').Value = '  -18.76%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.45'
$ws.Range('E37').Value = '  -22.46%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0353'
$ws.Range('E38').Value = '  -17.50%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.107'
$ws.Range('E39').Value = '  -12.81%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0695'
$ws.Range('E40').Value = '  -19.40%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.36'
$ws.Range('E41').Value = '  -20.59%  '

$ws.Range('B42').Value = 'USDe'
$ws.Range('C42').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  -0.03%  '

$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.542.74'
$ws.Range('E43').Value = '  -13.27%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.14'
$ws.Range('E44').Value = '  -22.47%  '

$ws.Range('B45').Value = 'Monero'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '107.18'
$ws.Range('E45').Value = '  -9.67%  '

$ws.Range('B46').Value = 'TheGraph'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.206'
$ws.Range('E46').Value = '  -21.79%  '

$ws.Range('B47').Value = 'BitgetToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.20'
$ws.Range('E47').Value = '  -4.44%  '

$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0958'
$ws.Range('E48').Value = '  -16.42%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.69'
$ws.Range('E49').Value = '  -21.20%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.50'
$ws.Range('E50').Value = '  -22.50%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.82'
$ws.Range('E51').Value = '  -25.64%  '
